# Apply the "Add 32 bit float for clarity" edit described by the diff:
#  1. Slide 4, shape "Rectangle 8": "How float value is represented " ->
#     "How 32 bit float value is represented "
#  2. Slide 5, shape "TextBox 6", last paragraph: merge the two runs
#     ("	- "" empty string " + "for strings") into a single run reading
#     "	- "" empty string for strings" and drop the now-redundant
#     trailing endParaRPr-only split.

$p = $ppt.ActivePresentation

# --- Change 1: slide 4 ----------------------------------------------------
$slide4 = $p.Slides.Item(4)
$rectShape = $slide4.Shapes.Item(5)
$rectShape.TextFrame.TextRange.Text = "How 32 bit float value is represented "

# --- Change 2: slide 5 ----------------------------------------------------
$slide5 = $p.Slides.Item(5)
$textBoxShape = $slide5.Shapes.Item(3)
$tr = $textBoxShape.TextFrame.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$fullRange = $lastPara.Characters(1, $lastPara.Length)
$fullRange.Text = "`t- " + [char]0x201C + [char]0x201D + " empty string for strings"
